# Turns the single run "Perfect. Below is a " into the proof-read markup
# Word produces once a typo "f" is typed in front of "Perfect" (making the
# misspelled "fPerfect", flagged by the spell-checker with <w:proofErr/>
# start/end tags) and the run subsequently gets split on the word boundary
# Word's live spell-check re-validates:
#
#   <w:proofErr w:type="spellStart"/>
#   <w:r>...<w:t>f</w:t></w:r>
#   <w:r>...<w:t>Perfect</w:t></w:r>
#   <w:proofErr w:type="spellEnd"/>
#   <w:r>...<w:t xml:space="preserve">. Below is a </w:t></w:r>
#
# Plain Range.Text / Find.Execute replacements only ever touch run *text* -
# they can't splice in bare <w:proofErr/> markers between runs, and any
# formatting-toggle trick used to force a run split gets coalesced straight
# back together on save once the toggle is undone. So instead we pull the
# paragraph's live OOXML out of the document, surgically rewrite just the
# one run, and hand the whole paragraph back via Range.InsertXML - which
# does let us place raw markup (including <w:proofErr/>) exactly where we
# want it.

function Find-LastIndexBefore($haystack, $needle, $before) {
    # .NET-style LastIndexOf(needle, startIndex) isn't reliably bounded in
    # this host, so walk forward matches ourselves and keep the last one
    # that still precedes $before.
    $cursor = -1
    $found = -1
    while ($true) {
        $cursor = $haystack.IndexOf($needle, $cursor + 1)
        if ($cursor -lt 0) { break }
        if ($cursor -ge $before) { break }
        $found = $cursor
    }
    return $found
}

$d = $word.ActiveDocument

# Round-trip through WordOpenXML to get the live serialized word/document.xml
# rather than hand-copying surrounding markup into this script.
$pkg = $d.Content.WordOpenXML
$partMarker = '<pkg:part pkg:name="/word/document.xml"'
$partIdx = $pkg.IndexOf($partMarker)
if ($partIdx -lt 0) { throw "word/document.xml part not found in package" }
$dataMarker = '<pkg:xmlData>'
$dataIdx = $pkg.IndexOf($dataMarker, $partIdx)
$docStart = $dataIdx + $dataMarker.Length
$dataEndMarker = '</pkg:xmlData></pkg:part>'
$docEnd = $pkg.IndexOf($dataEndMarker, $docStart)
$docXml = $pkg.Substring($docStart, $docEnd - $docStart)

# Isolate the paragraph that currently opens with "Perfect. Below is a ".
$paraMarker = '<w:p w14:paraId="58546E92"'
$pStart = $docXml.IndexOf($paraMarker)
if ($pStart -lt 0) { throw "Could not locate the target paragraph" }
$pEndTag = '</w:p>'
$pEndTagIdx = $docXml.IndexOf($pEndTag, $pStart)
$pEnd = $pEndTagIdx + $pEndTag.Length
$paraXml = $docXml.Substring($pStart, $pEnd - $pStart)

# Locate the run whose text is exactly "Perfect. Below is a ".
$targetText = 'Perfect. Below is a '
$tNeedle = '>' + $targetText + '</w:t>'
$tIdx = $paraXml.IndexOf($tNeedle)
if ($tIdx -lt 0) { throw "Could not locate run containing '$targetText'" }

$runStart = Find-LastIndexBefore $paraXml '<w:r>' $tIdx
$runStartAttr = Find-LastIndexBefore $paraXml '<w:r ' $tIdx
if ($runStartAttr -gt $runStart) { $runStart = $runStartAttr }
if ($runStart -lt 0) { throw "Could not locate start of target run" }

$runEndTag = '</w:r>'
$runEndTagIdx = $paraXml.IndexOf($runEndTag, $tIdx)
$runEnd = $runEndTagIdx + $runEndTag.Length
$oldRun = $paraXml.Substring($runStart, $runEnd - $runStart)

# Reuse the exact opening <w:r ...> tag and <w:rPr>...</w:rPr> block from the
# original run so every new run keeps identical formatting.
$rOpenTag = $oldRun.Substring(0, $oldRun.IndexOf('>') + 1)
$rPrStart = $oldRun.IndexOf('<w:rPr>')
$rPrEndIdx = $oldRun.IndexOf('</w:rPr>')
$rPrEnd = $rPrEndIdx + '</w:rPr>'.Length
$rPr = $oldRun.Substring($rPrStart, $rPrEnd - $rPrStart)

$newRuns = '<w:proofErr w:type="spellStart"/>' + `
    $rOpenTag + $rPr + '<w:t>f</w:t></w:r>' + `
    $rOpenTag + $rPr + '<w:t>Perfect</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    $rOpenTag + $rPr + '<w:t xml:space="preserve">. Below is a </w:t></w:r>'

$newParaXml = $paraXml.Substring(0, $runStart) + $newRuns + $paraXml.Substring($runEnd)

# InsertXML needs the fragment to carry its own namespace declarations since
# it's parsed standalone.
$nsDecl = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$newParaXmlNs = $newParaXml.Substring(0, 4) + $nsDecl + $newParaXml.Substring(4)

# Replace exactly this paragraph's text range (not the paragraph mark) with
# the rebuilt content.
$target = $d.Paragraphs(1).Range
$target.InsertXML($newParaXmlNs)

Write-Output $d.Paragraphs(1).Range.Text
